# Reposition / resize the floating text box ("Cuadro de texto 2") that
# sits behind the document text. The commit nudges it slightly left/down
# and widens it (height stays the same).
$d = $word.ActiveDocument

$shp = $d.Shapes(1)

# Values are expressed in points, matching the EMU offsets in the XML
# (1 pt = 12700 EMU):
#   left  : -148590 EMU (-11.7pt)  -> -150495 EMU (-11.85pt)
#   top   :  210820 EMU (16.6pt)   ->  212725 EMU (16.75pt)
#   width : 5703570 EMU (449.1pt)  -> 6065520 EMU (477.6pt)
#   height: 3842385 EMU (302.55pt) -> unchanged
$shp.Left = -11.85
$shp.Top = 16.75
$shp.Width = 477.6
